$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing values and add new ones in F2:H2
$ws.Range("B2").Value = 10.0
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 9
$ws.Range("H2").Value = 0

# Row 3: update B3 and clear C3:E3
$ws.Range("B3").Value = 7.0
$ws.Range("C3:E3").ClearContents()

# Row 4: update B4 and clear C4:E4
$ws.Range("B4").Value = 8.0
$ws.Range("C4:E4").ClearContents()

# Row 5: update B5 and clear C5:E5
$ws.Range("B5").Value = 6.0
$ws.Range("C5:E5").ClearContents()
